# tarificator : choix automatique du TYPM de photo / probleme famille du au fait
# qu'elle était au format numerique et 3 char
#
# Row 25 (HIKVISION FRANCE / HIK prefix) gains:
#   - D25 : "12s"  (a "temps" value, like the other rows in column D)
#   - F25 : "LA FAMILLE NE DOIT PAS ETRE EN NUMERIQUE " (a COMMENTAIRE note,
#           same family as the existing F-column warning about the CV4 family)
#
# Adding this text lengthens the longest strings in several columns, so the
# workbook's auto-fit column widths were refreshed as a side effect, and the
# active selection moved to F9:F10 (the next review target).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content edits -------------------------------------------------------
$ws.Range("D25").Value = "12s"
$ws.Range("F25").Value = "LA FAMILLE NE DOIT PAS ETRE EN NUMERIQUE "

# --- column width refresh (auto-fit after the longer text was entered) ---
$ws.Columns.Item(1).ColumnWidth = 42.166666666666664
$ws.Columns.Item(2).ColumnWidth = 34.833333333333336
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(5).ColumnWidth = 16
$ws.Columns.Item(6).ColumnWidth = 47
$ws.Columns.Item(7).ColumnWidth = 31.333333333333332

# --- selection moved on to the next rows to review ------------------------
$ws.Range("F9:F10").Select()
